$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ASML)
$ws.Range("K2").Value = 59.6
$ws.Range("N2").Value = 54.77309453746771

# Row 3 (TSM)
$ws.Range("K3").Value = 58.4
$ws.Range("N3").Value = 54.77309453746771

# Row 4 (QCOM)
$ws.Range("E4").Value = 51
$ws.Range("F4").Value = 6.12
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 43
$ws.Range("K4").Value = 51.6
$ws.Range("N4").Value = 54.77309453746771

# Row 5 (NVDA)
$ws.Range("F5").Value = 1.74
$ws.Range("K5").Value = 49.4
$ws.Range("N5").Value = 54.77309453746771

# Row 6 (AMD)
$ws.Range("K6").Value = 46.4
$ws.Range("N6").Value = 54.77309453746771
